$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the hours worked (End - Start) for the row 7 entry and its description
$ws.Range("C7").Value = 15
$ws.Range("E7").Value = "Routing, styling, invoking classes. Updates found at https://github.com/leono93/flutter-project"

# Update the active selection to match the author's cursor position
$ws.Range("E10").Select()
